# amounttomove and thresholdtomove are now a part of the input.csv system
#
# Insert two new columns (M:N) into the sheet, shifting the existing
# "numbspecies" column (and everything after it) two columns to the right.
# Populate the new columns with header labels + per-row values, matching
# the style of the column that used to live at M (style index 2 / yellow fill).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns before the old column M ("numbspecies"),
# pushing M:AW to O:AY.
$ws.Range("M1:N7").Insert(-4161)

# Header row: new labels for the inserted columns.
$ws.Range("M1").Value = "thresholdtomove"
$ws.Range("N1").Value = "amounttomove"

# Match the highlighted style used by the neighbouring "numbspecies" column.
$ws.Range("M1:N7").Style = $ws.Range("O1").Style

# Data rows: thresholdtomove = 30, amounttomove = 0.25 for every data row.
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 13).Value = 30
    $ws.Cells.Item($row, 14).Value = 0.25
}

# Keep the selection where the author left it after the edit.
$ws.Range("M2:N7").Select()
